$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1343.6666
$ws.Range("J107").Value = 1463.8889
$ws.Range("L107").Value = 1463.8889
$ws.Range("N107").Value = -5303.8889

$ws.Range("H116").Value = 2759.92
$ws.Range("I116").Value = 2704.041
$ws.Range("J116").Value = 4799.5
$ws.Range("K116").Value = 2704.041
$ws.Range("L116").Value = 4799.5
$ws.Range("M116").Value = 737.9589999999998
$ws.Range("N116").Value = -11683.5

$ws.Range("H138").Value = 2487.3333
$ws.Range("I138").Value = 2993
$ws.Range("K138").Value = 8979
$ws.Range("M138").Value = -3839

$ws.Range("H141").Value = 4607
$ws.Range("I141").Value = 3913
$ws.Range("K141").Value = 11739
$ws.Range("M141").Value = -6559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N23").ClearContents()
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -4741

$ws.Range("H32").Value = 15039.895
$ws.Range("I32").Value = 5181
$ws.Range("J32").Value = 31940.857
$ws.Range("K32").Value = 5181
$ws.Range("L32").Value = 31940.857
$ws.Range("M32").Value = -4894
$ws.Range("N32").Value = -32514.857

$ws.Range("H97").Value = 3562.111
$ws.Range("I97").Value = 1867.15
$ws.Range("J97").Value = 8404.857
$ws.Range("K97").Value = 1867.15
$ws.Range("L97").Value = 8404.857
$ws.Range("M97").Value = -1371.15
$ws.Range("N97").Value = -9396.857

$ws.Range("H109").Value = 295000
$ws.Range("J109").Value = 295000
$ws.Range("L109").Value = 295000
$ws.Range("N109").Value = -297774

$ws.Range("H137").Value = 39983.688
$ws.Range("J137").Value = 42973.9
$ws.Range("L137").Value = 42973.9
$ws.Range("N137").Value = -53173.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1950
$ws.Range("I5").Value = 1900
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1900
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -1787
$ws.Range("N5").Value = -2226

$ws.Range("H96").Value = 43333.168
$ws.Range("I96").Value = 10000
$ws.Range("K96").Value = 10000
$ws.Range("M96").Value = -7254

$ws.Range("H99").Value = 1695
$ws.Range("I99").Value = 1180.909
$ws.Range("K99").Value = 1180.909
$ws.Range("M99").Value = 317.0909999999999

$ws.Range("H134").Value = 15964
$ws.Range("I134").Value = 9070.532999999999
$ws.Range("K134").Value = 27211.599
$ws.Range("M134").Value = -24676.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 358598.34
$ws.Range("J19").Value = 737.5
$ws.Range("L19").Value = 737.5
$ws.Range("N19").Value = -1077.5

$ws.Range("H24").Value = 358598.34
$ws.Range("J24").Value = 737.5
$ws.Range("L24").Value = 737.5
$ws.Range("N24").Value = -1077.5

$ws.Range("N44").ClearContents()
$ws.Range("H44").Value = 27500
$ws.Range("I44").Value = 27500
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 27500
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -27058

$ws.Range("H86").Value = 6463
$ws.Range("I86").Value = 5606.857
$ws.Range("K86").Value = 5606.857
$ws.Range("M86").Value = -4483.857

$ws.Range("H89").Value = 6463
$ws.Range("I89").Value = 5606.857
$ws.Range("K89").Value = 28034.285
$ws.Range("M89").Value = -22418.285

$ws.Range("H122").Value = 4667.1724
$ws.Range("J122").Value = 11692.571
$ws.Range("L122").Value = 35077.713
$ws.Range("N122").Value = -39977.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 729.5
$ws.Range("J5").Value = 6102897
$ws.Range("K5").Value = 2188.5
$ws.Range("L5").Value = 18308691
$ws.Range("M5").Value = -2076.5
$ws.Range("N5").Value = -18308915

$ws.Range("H17").Value = 91
$ws.Range("J17").Value = 91
$ws.Range("L17").Value = 273
$ws.Range("N17").Value = -611

$ws.Range("H22").Value = 7999.5
$ws.Range("J22").Value = 6000
$ws.Range("L22").Value = 18000
$ws.Range("N22").Value = -18338

$ws.Range("H27").Value = 7999.5
$ws.Range("J27").Value = 6000
$ws.Range("L27").Value = 18000
$ws.Range("N27").Value = -18204

$ws.Range("H68").Value = 15360.5
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5189

$ws.Range("H71").Value = 15360.5
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 18000
$ws.Range("M71").Value = -13944

$ws.Range("H75").Value = 224.5
$ws.Range("J75").Value = 199.33333
$ws.Range("L75").Value = 597.99999
$ws.Range("N75").Value = -2593.99999

$ws.Range("H78").Value = 224.5
$ws.Range("J78").Value = 199.33333
$ws.Range("L78").Value = 1793.99997
$ws.Range("N78").Value = -11777.99997

$ws.Range("H88").Value = 28428.428
$ws.Range("J88").Value = 28428.428
$ws.Range("L88").Value = 85285.284
$ws.Range("N88").Value = -86141.284

$ws.Range("H91").Value = 28428.428
$ws.Range("J91").Value = 28428.428
$ws.Range("L91").Value = 85285.284
$ws.Range("N91").Value = -88249.284

$ws.Range("H107").Value = 2841803.5
$ws.Range("I107").Value = 590.3333
$ws.Range("J107").Value = 6251259.5
$ws.Range("K107").Value = 1770.9999
$ws.Range("L107").Value = 18753778.5
$ws.Range("M107").Value = 149.0001
$ws.Range("N107").Value = -18757618.5

$ws.Range("M118").ClearContents()
$ws.Range("H118").Value = 15516
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 15516
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 46548
$ws.Range("N118").Value = -49034

$ws.Range("I135").Value = 729.5
$ws.Range("J135").Value = 6102897
$ws.Range("K135").Value = 6565.5
$ws.Range("L135").Value = 54926073
$ws.Range("M135").Value = -4030.5
$ws.Range("N135").Value = -54931143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30174.908
$ws.Range("J52").Value = 30174.908
$ws.Range("L52").Value = 30174.908
$ws.Range("N52").Value = -30692.908

$ws.Range("H102").Value = 7979.278
$ws.Range("I102").Value = 5683
$ws.Range("J102").Value = 26349.5
$ws.Range("K102").Value = 5683
$ws.Range("L102").Value = 26349.5
$ws.Range("M102").Value = -4061
$ws.Range("N102").Value = -29593.5

$ws.Range("H113").Value = 4436.778
$ws.Range("I113").Value = 3643.3845
$ws.Range("K113").Value = 3643.3845
$ws.Range("M113").Value = -1473.3845

$ws.Range("H122").Value = 5715.125
$ws.Range("I122").Value = 1912.1111
$ws.Range("K122").Value = 5736.3333
$ws.Range("M122").Value = -3286.3333

$ws.Range("H132").Value = 12511.25
$ws.Range("I132").Value = 7349.75
$ws.Range("J132").Value = 27995.75
$ws.Range("K132").Value = 22049.25
$ws.Range("L132").Value = 83987.25
$ws.Range("M132").Value = -19519.25
$ws.Range("N132").Value = -89047.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2861000
$ws.Range("I2").Value = 3400
$ws.Range("K2").Value = 3400
$ws.Range("M2").Value = -3288

$ws.Range("H23").Value = 3006666.8
$ws.Range("I23").Value = 3006666.8
$ws.Range("K23").Value = 3006666.8
$ws.Range("M23").Value = -3006436.8

$ws.Range("H40").Value = 7928.722
$ws.Range("I40").Value = 2851.5
$ws.Range("K40").Value = 2851.5
$ws.Range("M40").Value = -2715.5

$ws.Range("H82").Value = 5083.8076
$ws.Range("I82").Value = 2693.182
$ws.Range("J82").Value = 6836.933
$ws.Range("K82").Value = 2693.182
$ws.Range("L82").Value = 6836.933
$ws.Range("M82").Value = -2332.182
$ws.Range("N82").Value = -7558.933

$ws.Range("H85").Value = 5083.8076
$ws.Range("I85").Value = 2693.182
$ws.Range("J85").Value = 6836.933
$ws.Range("K85").Value = 2693.182
$ws.Range("L85").Value = 6836.933
$ws.Range("M85").Value = -1445.182
$ws.Range("N85").Value = -9332.933000000001

$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0

$ws.Range("H132").Value = 651804.7
$ws.Range("I132").Value = 2478.2368
$ws.Range("K132").Value = 7434.7104
$ws.Range("M132").Value = -4904.7104

$ws.Range("H136").Value = 12625.31
$ws.Range("I136").Value = 10278.962
$ws.Range("K136").Value = 30836.886
$ws.Range("M136").Value = -28286.886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 226507.47
$ws.Range("I2").Value = 2258
$ws.Range("K2").Value = 2258
$ws.Range("M2").Value = -2146

$ws.Range("H53").Value = 25000
$ws.Range("J53").Value = 35000
$ws.Range("L53").Value = 35000
$ws.Range("N53").Value = -36214

$ws.Range("H96").Value = 1821.4117
$ws.Range("I96").Value = 1379.4286
$ws.Range("J96").Value = 2130.8
$ws.Range("K96").Value = 1821.4117
$ws.Range("L96").Value = 2130.8
$ws.Range("M96").Value = -6.42859999999996
$ws.Range("N96").Value = -4876.8

